# Refactor EEG data processing scripts and add README.md for improved
# clarity and functionality of pipeline.
#
# Excel-side edit to matlab/subject_data_info.xlsx:
#  - widen the descriptive columns (A-F) on Sheet1 so the newly-clarified
#    column headers / longer notes text are readable
#  - leave the current selection on the populated data block (B2:F15)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (Excel stores widths in "characters" of the Normal
# style font; ColumnWidth here is that same character-width unit) ---
$ws.Columns.Item(1).ColumnWidth = 30.666666666666668   # A: Subject -> wider
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666   # B: ID -> wider
$ws.Columns.Item(3).ColumnWidth = 11.5                 # C: Session (unchanged width, now its own col)
$ws.Columns.Item(4).ColumnWidth = 13.166666666666666   # D: Run/condition -> wider
$ws.Columns.Item(5).ColumnWidth = 14.5                 # E: Type -> wider
$ws.Columns.Item(6).ColumnWidth = 94.5                 # F: Notes -> much wider

# --- Selection: active cell B2, selected range B2:F15 ---
$ws.Range("B2:F15").Select()
